{"js": "// Update intro paragraph: drop \"och fridlysta arter\" + trailing BILAGA 1 sentence.\nconst body = context.document.body;\nbody.paragraphs.load(\"items,text,style\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\nconst OLD_INTRO =\n  \"Nedan presenteras fynd av naturv\u00e5rdsarter och fridlysta arter som gjorts i det avverkningsanm\u00e4lda omr\u00e5det, samt relevanta utdrag ur standarderna f\u00f6r FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter.\";\nconst NEW_INTRO =\n  \"Nedan presenteras fynd av naturv\u00e5rdsarter som gjorts i det avverkningsanm\u00e4lda omr\u00e5det, samt relevanta utdrag ur standarderna f\u00f6r FSC, Chain of Custody, Controlled Wood och PEFC.\";\n\nconst toDelete = [];\nfor (let i = 0; i < items.length; i++) {\n  const p = items[i];\n  const text = p.text;\n\n  if (text === OLD_INTRO) {\n    p.insertText(NEW_INTRO, Word.InsertLocation.replace);\n  } else if (\n    // Entire \"Fridlysta arter\" sub-section (heading through the kn\u00e4rot\n    // figure caption), now superseded by its own appendix document.\n    text === \"Fridlysta arter\" ||\n    text ===\n      \"F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er och v\u00e4xtplatser i den avverkningsanm\u00e4lda skogen: \" ||\n    text ===\n      \"I det avverkningsanm\u00e4lda omr\u00e5det finns 0 fyndplatser f\u00f6r kn\u00e4rot registrerade p\u00e5 Artportalen. F\u00f6r att kn\u00e4rotslokaler inte ska f\u00f6rsvinna vid avverkning kr\u00e4vs att en buffertzon p\u00e5 50 m l\u00e4mnas kring samtliga fyndplatser. Figur 2 visar gr\u00e4nserna f\u00f6r dessa buffertzoner.\" ||\n    text ===\n      \"Figur 2. Fyndplatser och buffertzoner f\u00f6r kn\u00e4rot i det avverkningsanm\u00e4lda omr\u00e5det. Endast fyndplaster vars buffertzoner \u00f6verlappar med det avverkningsanm\u00e4lda omr\u00e5det har tagits med i visualiseringen. Kartans mittpunktskoordinat \u00e4r N\\u00A06675710, E\\u00A0656144 i SWEREF 99 TM.\" ||\n    // Trailing page break + \"BILAGA 1\" title paragraph.\n    text === \"\\f\" ||\n    text === \"BILAGA 1 - Fridlysta arter\"\n  ) {\n    toDelete.push(p);\n  } else if (text === \"\" && p.style === \"Caption\") {\n    // The kn\u00e4rot figure's (Picture 2) own caption-styled paragraph \u2014 empty\n    // text because the picture lives in an inline drawing, not run text.\n    // Only the second (kn\u00e4rot) figure's image paragraph should go; the\n    // first (naturv\u00e5rdsarter) figure is kept. Disambiguate via the\n    // paragraph that follows.\n    const next = items[i + 1];\n    if (\n      next &&\n      next.text ===\n        \"Figur 2. Fyndplatser och buffertzoner f\u00f6r kn\u00e4rot i det avverkningsanm\u00e4lda omr\u00e5det. Endast fyndplaster vars buffertzoner \u00f6verlappar med det avverkningsanm\u00e4lda omr\u00e5det har tagits med i visualiseringen. Kartans mittpunktskoordinat \u00e4r N\\u00A06675710, E\\u00A0656144 i SWEREF 99 TM.\"\n    ) {\n      toDelete.push(p);\n    }\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\n// The \"first page\" header carries the document date in its own body.\nconst firstHeader = context.document.sections.getFirst().getHeader(Word.HeaderFooterType.firstPage);\nfirstHeader.search(\"2023-09-06\", { matchCase: true }).load(\"items,text\");\nawait context.sync();\n\nfirstHeader.search(\"2023-09-06\", { matchCase: true }).items.forEach((r) => {\n  r.insertText(\"2023-09-08\", Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) equivalent of edit.js.\n$d = $word.ActiveDocument\n\n# The figure captions use a non-breaking space (U+00A0) before the\n# coordinate digits (e.g. \"N<nbsp>6675710\"); build it explicitly so the\n# literal below matches the run text exactly.\n$nbsp = [char]0x00A0\n\n$oldIntro = \"Nedan presenteras fynd av naturv\u00e5rdsarter och fridlysta arter som gjorts i det avverkningsanm\u00e4lda omr\u00e5det, samt relevanta utdrag ur standarderna f\u00f6r FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter.\"\n$newIntro = \"Nedan presenteras fynd av naturv\u00e5rdsarter som gjorts i det avverkningsanm\u00e4lda omr\u00e5det, samt relevanta utdrag ur standarderna f\u00f6r FSC, Chain of Custody, Controlled Wood och PEFC.\"\n\n$knarotCaption = \"Figur 2. Fyndplatser och buffertzoner f\u00f6r kn\u00e4rot i det avverkningsanm\u00e4lda omr\u00e5det. Endast fyndplaster vars buffertzoner \u00f6verlappar med det avverkningsanm\u00e4lda omr\u00e5det har tagits med i visualiseringen. Kartans mittpunktskoordinat \u00e4r N\" + $nbsp + \"6675710, E\" + $nbsp + \"656144 i SWEREF 99 TM.\"\n\n# Exact paragraph texts (paragraph mark stripped) that make up the whole\n# \"Fridlysta arter\" sub-section, plus the trailing page-break + \"BILAGA 1\"\n# title paragraph at the very end of the body \u2014 all dropped in this edit.\n$deleteTexts = @(\n  \"Fridlysta arter\",\n  \"F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er och v\u00e4xtplatser i den avverkningsanm\u00e4lda skogen: \",\n  \"I det avverkningsanm\u00e4lda omr\u00e5det finns 0 fyndplatser f\u00f6r kn\u00e4rot registrerade p\u00e5 Artportalen. F\u00f6r att kn\u00e4rotslokaler inte ska f\u00f6rsvinna vid avverkning kr\u00e4vs att en buffertzon p\u00e5 50 m l\u00e4mnas kring samtliga fyndplatser. Figur 2 visar gr\u00e4nserna f\u00f6r dessa buffertzoner.\",\n  $knarotCaption,\n  [string][char]12,\n  \"BILAGA 1 - Fridlysta arter\"\n)\n\n$count = $d.Paragraphs.Count\n$indicesToDelete = New-Object System.Collections.ArrayList\n\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $text = $p.Range.Text.TrimEnd(\"`r\")\n\n  if ($text -eq $oldIntro) {\n    $p.Range.Text = $newIntro\n  } elseif ($deleteTexts -contains $text) {\n    [void]$indicesToDelete.Add($i)\n  } elseif ($text -eq \"\" -and $p.Range.ParagraphStyle.NameLocal -eq \"Caption\") {\n    # The kn\u00e4rot figure's own (picture-only) caption-styled paragraph has no\n    # run text \u2014 disambiguate it from the naturv\u00e5rdsarter figure's picture\n    # paragraph by checking whether the very next paragraph is its caption.\n    if ($i -lt $count) {\n      $nextText = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd(\"`r\")\n      if ($nextText -eq $knarotCaption) {\n        [void]$indicesToDelete.Add($i)\n      }\n    }\n  }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sorted = $indicesToDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n  $d.Paragraphs.Item($idx).Range.Delete()\n}\n\n# The \"first page\" header carries the document date in its own story\n# (wdHeaderFooterFirstPage = 2; the document has titlePg / a distinct\n# first-page header, which is where \"2023-09-06\" actually lives).\n$firstHeader = $d.Sections.Item(1).Headers.Item(2)\n$rng = $firstHeader.Range\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"2023-09-06\"\n$rng.Find.Replacement.Text = \"2023-09-08\"\n$rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 0, $false, $null, 2) | Out-Null\n"}
